# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column between
# "In Advance" (M) and "Late" (N): everything from N onward shifts one
# column to the right, and the sheet becomes the active tab/selection
# (previously "NewLoanInput" was the active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting
# N/O/P -> O/P/Q. Cell values/styles to the right of the insertion point
# move along with their columns; the newly created column N is empty.
$ws.Columns("N").Insert() | Out-Null

# Match the width of the new column to its left-hand neighbour (column M),
# same as the width Excel carries over for an inserted column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with R6 selected
# (previously NewLoanInput was active, and R7 was selected on this sheet).
$ws.Activate()
$ws.Range("R6").Select() | Out-Null
